# Update the FRED RRPONTSYD workbook:
#  - Append new daily observations (rows 519-531) to the "Data" sheet
#  - Refresh the FRED series metadata on the "SeriesInfo" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Data" sheet - append the new observations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

$lastRow = 518

# Give the new date cells (column A) the same style/number-format as the
# existing date column before filling in values.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + ($lastRow + 1) + ":A" + ($lastRow + 13)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(45229, 1138.035),
    @(45230, 1137.697),
    @(45231, 1079.462),
    @(45232, 1054.986),
    @(45233, 1071.139),
    @(45236, 1062.878),
    @(45237, 1008.685),
    @(45238, 1024.451),
    @(45239, 993.314),
    @(45240, 1032.72),
    @(45243, 1020.272),
    @(45244, 988.298),
    @(45245, 944.241)
)

$r = $lastRow + 1
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. "SeriesInfo" sheet - refresh metadata values
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SeriesInfo")

# Use an out-of-the-way scratch cell to stage plain-text values (prefixed
# with an apostrophe) and paste just the values in, so the date-like
# strings land as literal text (matching the existing inlineStr cells)
# instead of being auto-converted to date serials, and without picking
# up any new cell style.
$helper = $ws2.Cells.Item(100, 1)

$helper.Value = "'2023-11-15"
$helper.Copy()
$ws2.Cells.Item(3, 2).PasteSpecial(-4163)
$ws2.Cells.Item(4, 2).PasteSpecial(-4163)
$ws2.Cells.Item(7, 2).PasteSpecial(-4163)

$helper.Value = "'2023-11-15 13:01:02-06"
$helper.Copy()
$ws2.Cells.Item(14, 2).PasteSpecial(-4163)

$excel.CutCopyMode = 0
$helper.Clear()

$ws2.Cells.Item(15, 2).Value = 93
